$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (D) and volume-change (E) values for the cryptos table.
# Column D cells that contain purely numeric-looking text are forced to
# remain Text so Excel does not auto-convert them to numbers and lose
# the exact displayed formatting (e.g. trailing zeros, dotted thousands).

$ws.Range("D2").Value = "29.252.31"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.903.23"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.20"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4649"
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3918"
$ws.Range("E8").Value = "  -0.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07898"
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.85"
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("D12").Value = "1.905.69"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.084"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.751"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06995"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.33"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009987"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.12"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D21").Value = "29.253.09"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.305"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "2.139.23"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.21"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.986"
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "118.90"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  -5.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09334"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9034"
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.263"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.328"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.206"
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.183"
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05777"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.727"
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5710"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.727"
$ws.Range("E43").Value = "  -2.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.98"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5364"
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.182"
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07024"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.576"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.34"
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.061"
$ws.Range("E51").Value = "  +0.16%  "
